# Daily attendance processing - 2026-01-04 14:03:33
# Rotate the "Recorded By" (column G) list of names/emails for each data row:
# move the first comma-separated entry to the end of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value()
    if ($val -ne $null) {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $newVal = $rotated -join ", "
            $cell.Value = $newVal
        }
    }
}
